# edit.ps1 - applies the changes described by the diff:
#  1. Fixes a typo in a data string on the "sample01" sheet (cell I16):
#       "...238,2239,240..." -> "...238,223,240..."
#  2. Adds a new worksheet "sample02" (placed after "sample01") that
#     shows a second example message, built from the same layout as
#     "sample01" but with a few different rows/values.
#  3. Restores view-state (selection/active sheet) as closely as the
#     object model allows.

$wb = $excel.ActiveWorkbook

$wsSchema   = $wb.Worksheets.Item(1)   # messageSchema
$wsSample01 = $wb.Worksheets.Item(2)   # sample01

# -----------------------------------------------------------------
# 1. Fix the typo inside the existing sample01 sheet
# -----------------------------------------------------------------
$wsSample01.Range("I16").Value = "[303,132,106,347,236,237,238,223,240,241]"

# -----------------------------------------------------------------
# 2. Build the new "sample02" worksheet, positioned right after
#    "sample01"
# -----------------------------------------------------------------
$wsSample02 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wsSample01)
$wsSample02.Name = "sample02"

# Copy over the values from sample01 first (rows 1-10 share the same
# layout), then adjust the handful of cells that differ.
$wsSample01.Range("A1:J10").Copy()
$wsSample02.Range("A1:J10").PasteSpecial(-4104) | Out-Null   # xlPasteAll

# Copy the column formatting used for I:J (style for the
# "message"/"replay" example columns) down through row 17 as well,
# so the empty placeholder rows keep the correct look.
$wsSample01.Range("I1:J17").Copy()
$wsSample02.Range("I1:J17").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$excel.CutCopyMode = 0

# Row 10 differs from sample01: the key value is 10000 (not 20000)
$wsSample02.Range("A10").Value = 10000

# Rows 11-17 in sample02 are blank placeholder rows (only the I/J
# formatting carries over, already applied above).
$wsSample02.Range("A11:H17").ClearContents()
$wsSample02.Range("I11:J16").ClearContents()

# -----------------------------------------------------------------
# 3. View-state: selection / active sheet for each tab
# -----------------------------------------------------------------
$wsSchema.Activate()
$wsSchema.Range("A8").Select()

$wsSample01.Activate()
$wsSample01.Range("A1:J17").Select()

$wsSample02.Activate()
$wsSample02.Range("A10").Select()

$wsSchema.Activate()
